$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (matches original inlineStr semantics),
# without leaving a residual NumberFormat style on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "24.100.84"
Set-TextValue $ws.Range("E2") "  -3.70%  "

Set-TextValue $ws.Range("D3") "1.645.55"
Set-TextValue $ws.Range("E3") "  -3.51%  "

Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  -0.05%  "

Set-TextValue $ws.Range("D5") "307.54"
Set-TextValue $ws.Range("E5") "  -2.86%  "

Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.02%  "

Set-TextValue $ws.Range("D7") "0.3909"
Set-TextValue $ws.Range("E7") "  -2.09%  "

Set-TextValue $ws.Range("D8") "0.3856"
Set-TextValue $ws.Range("E8") "  -4.68%  "

Set-TextValue $ws.Range("D9") "1.001"
Set-TextValue $ws.Range("E9") "  -0.11%  "

Set-TextValue $ws.Range("D10") "1.349"
Set-TextValue $ws.Range("E10") "  -8.47%  "

Set-TextValue $ws.Range("D11") "49.08"
Set-TextValue $ws.Range("E11") "  -7.36%  "

Set-TextValue $ws.Range("D12") "0.08468"
Set-TextValue $ws.Range("E12") "  -3.96%  "

Set-TextValue $ws.Range("D13") "23.85"
Set-TextValue $ws.Range("E13") "  -8.52%  "

Set-TextValue $ws.Range("D14") "7.128"
Set-TextValue $ws.Range("E14") "  -4.82%  "

Set-TextValue $ws.Range("D15") "0.00001283"
Set-TextValue $ws.Range("E15") "  -4.98%  "

Set-TextValue $ws.Range("D16") "7.489"
Set-TextValue $ws.Range("E16") "  -6.20%  "

Set-TextValue $ws.Range("D17") "1.644.44"
Set-TextValue $ws.Range("E17") "  -2.64%  "

Set-TextValue $ws.Range("D18") "94.88"
Set-TextValue $ws.Range("E18") "  -1.22%  "

Set-TextValue $ws.Range("D19") "0.06969"

Set-TextValue $ws.Range("D20") "20.74"
Set-TextValue $ws.Range("E20") "  -0.27%  "

Set-TextValue $ws.Range("E22") "  -0.03%  "

Set-TextValue $ws.Range("D23") "13.63"
Set-TextValue $ws.Range("E23") "  -4.83%  "

Set-TextValue $ws.Range("D24") "24.096.18"
Set-TextValue $ws.Range("E24") "  -3.75%  "

Set-TextValue $ws.Range("D25") "2.327"
Set-TextValue $ws.Range("E25") "  -3.11%  "

Set-TextValue $ws.Range("D26") "2.705"
Set-TextValue $ws.Range("E26") "  -8.03%  "

Set-TextValue $ws.Range("D27") "22.46"
Set-TextValue $ws.Range("E27") "  -4.83%  "

Set-TextValue $ws.Range("D28") "158.08"
Set-TextValue $ws.Range("E28") "  -2.95%  "

Set-TextValue $ws.Range("D29") "8.753"
Set-TextValue $ws.Range("E29") "  +3.61%  "

Set-TextValue $ws.Range("D30") "141.47"
Set-TextValue $ws.Range("E30") "  -7.05%  "

Set-TextValue $ws.Range("D31") "5.297"
Set-TextValue $ws.Range("E31") "  -12.44%  "

Set-TextValue $ws.Range("D32") "2.461"
Set-TextValue $ws.Range("E32") "  -8.36%  "

Set-TextValue $ws.Range("D33") "1.827.48"
Set-TextValue $ws.Range("E33") "  -4.02%  "

Set-TextValue $ws.Range("D34") "6.906"
Set-TextValue $ws.Range("E34") "  -4.17%  "

Set-TextValue $ws.Range("D35") "0.08020"
Set-TextValue $ws.Range("E35") "  -7.18%  "

Set-TextValue $ws.Range("D36") "0.02913"
Set-TextValue $ws.Range("E36") "  -7.99%  "

Set-TextValue $ws.Range("B37") "Algorand"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D37") "0.2698"
Set-TextValue $ws.Range("E37") "  -7.68%  "

Set-TextValue $ws.Range("B38") "ImmutableX"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.9577"
Set-TextValue $ws.Range("E38") "  -8.75%  "

Set-TextValue $ws.Range("D39") "0.09203"
Set-TextValue $ws.Range("E39") "  -5.26%  "

Set-TextValue $ws.Range("D40") "1.460"
Set-TextValue $ws.Range("E40") "  -1.48%  "

Set-TextValue $ws.Range("D41") "9.962"
Set-TextValue $ws.Range("E41") "  -9.73%  "

Set-TextValue $ws.Range("D42") "0.7609"
Set-TextValue $ws.Range("E42") "  -8.00%  "

Set-TextValue $ws.Range("E43") "  -6.91%  "

Set-TextValue $ws.Range("D44") "16.11"
Set-TextValue $ws.Range("E44") "  -5.40%  "

Set-TextValue $ws.Range("D45") "0.6915"
Set-TextValue $ws.Range("E45") "  -6.27%  "

Set-TextValue $ws.Range("D46") "2.482"
Set-TextValue $ws.Range("E46") "  -7.74%  "

Set-TextValue $ws.Range("D47") "4.101"
Set-TextValue $ws.Range("E47") "  -3.49%  "

Set-TextValue $ws.Range("E48") "  +0.02%  "

Set-TextValue $ws.Range("D49") "0.08349"
Set-TextValue $ws.Range("E49") "  -9.54%  "

Set-TextValue $ws.Range("D50") "133.90"
Set-TextValue $ws.Range("E50") "  -4.41%  "

Set-TextValue $ws.Range("D51") "1.262"
Set-TextValue $ws.Range("E51") "  -10.15%  "
